$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.204.97"
$ws.Range("E2").Value = "  -3.11%  "
$ws.Range("D3").Value = "1.648.04"
$ws.Range("E3").Value = "  -3.71%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.38%  "
$ws.Range("D5").Value = "'308.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.92%  "
$ws.Range("D6").Value = "'1.003"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.23%  "
$ws.Range("D7").Value = "'0.3880"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.06%  "
$ws.Range("D8").Value = "'0.3876"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.00%  "
$ws.Range("D9").Value = "'1.004"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.42%  "
$ws.Range("D10").Value = "'1.365"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -9.02%  "
$ws.Range("D11").Value = "'49.22"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.72%  "
$ws.Range("D12").Value = "'0.08480"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.06%  "
$ws.Range("D13").Value = "'24.42"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -6.41%  "
$ws.Range("D14").Value = "'7.156"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.37%  "
$ws.Range("D15").Value = "'0.00001288"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.16%  "
$ws.Range("D16").Value = "'7.527"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.13%  "
$ws.Range("D17").Value = "1.642.09"
$ws.Range("E17").Value = "  -4.32%  "
$ws.Range("D18").Value = "'94.90"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.82%  "
$ws.Range("E19").Value = "  -4.14%  "
$ws.Range("D20").Value = "'20.86"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.85%  "
$ws.Range("D21").Value = "'6.964"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.04%  "
$ws.Range("D22").Value = "'1.004"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.52%  "
$ws.Range("E23").Value = "  -4.93%  "
$ws.Range("D24").Value = "24.200.80"
$ws.Range("E24").Value = "  -3.13%  "
$ws.Range("D25").Value = "'2.367"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.48%  "
$ws.Range("D26").Value = "'2.747"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -8.13%  "
$ws.Range("D27").Value = "'22.61"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.47%  "
$ws.Range("D28").Value = "'158.09"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.08%  "
$ws.Range("B29").Value = "BitcoinCash"
$ws.Range("C29").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D29").Value = "'142.43"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.37%  "
$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D30").Value = "'8.295"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.64%  "
$ws.Range("D31").Value = "'5.358"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -13.72%  "
$ws.Range("D32").Value = "'2.502"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.97%  "
$ws.Range("D33").Value = "1.823.36"
$ws.Range("E33").Value = "  -4.29%  "
$ws.Range("D34").Value = "'0.08138"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.86%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "'0.9921"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.13%  "
$ws.Range("B36").Value = "InternetComputer(DFINITY)"
$ws.Range("C36").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D36").Value = "'6.842"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.30%  "
$ws.Range("D37").Value = "'0.02945"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.89%  "
$ws.Range("D38").Value = "'0.2729"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.58%  "
$ws.Range("D39").Value = "'0.09312"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.27%  "
$ws.Range("D40").Value = "'1.473"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.69%  "
$ws.Range("D41").Value = "'9.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -7.74%  "
$ws.Range("D42").Value = "'0.7667"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.27%  "
$ws.Range("D43").Value = "'13.16"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.78%  "
$ws.Range("D44").Value = "'16.04"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.82%  "
$ws.Range("D45").Value = "'2.502"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.43%  "
$ws.Range("E46").Value = "  -6.73%  "
$ws.Range("D47").Value = "'4.097"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.91%  "
$ws.Range("D48").Value = "'1.002"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.13%  "
$ws.Range("D49").Value = "'0.08479"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.76%  "
$ws.Range("B50").Value = "Flow"
$ws.Range("C50").Value = "https://coinranking.com/coin/QQ0NCmjVq+flow-flow"
$ws.Range("D50").Value = "'1.267"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -8.31%  "
$ws.Range("B51").Value = "Quant"
$ws.Range("C51").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D51").Value = "'133.94"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.93%  "
